$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 161..245 (after the weekly shift): row|D(date)|J(volumen)|K(precio min)|L(precio max)|M(precio prom)|O(origen)|P(precio $/Kg)
$rowsData = @(
    "161|45134|2000|2400|2500|2450|Provincia del Elquí|1633",
    "162|44894|2500|1500|2000|1750|Provincia del Elquí|1167",
    "163|45050|2400|2000|2500|2250|Provincia del Elquí|1500",
    "164|44649|2800|2300|2500|2400|Provincia del Elquí|1600",
    "165|44754|2400|1500|2000|1750|Provincia del Elquí|1167",
    "166|44789|3000|2000|2500|2250|Provincia del Elquí|1500",
    "167|44299|3400|2000|2500|2250|Provincia del Elquí|1500",
    "168|44334|3440|1300|1500|1400|Provincia del Elquí|933",
    "169|44210|3200|1300|1500|1400|Provincia del Elquí|933",
    "170|44483|3300|1500|2000|1750|Provincia del Elquí|1167",
    "171|44875|2560|1500|2000|1750|Provincia del Elquí|1167",
    "172|44264|3600|2000|2500|2250|Provincia del Elquí|1500",
    "173|44306|3400|2000|2500|2250|Provincia del Elquí|1500",
    "174|44516|3100|1300|1500|1400|Provincia del Elquí|933",
    "175|44435|6560|2000|2500|2250|Provincia del Elquí|1500",
    "176|44588|3200|2500|3000|2750|Provincia del Elquí|1833",
    "177|45120|2000|2500|3000|2750|Provincia del Elquí|1833",
    "178|44236|3200|1500|2000|1750|Provincia del Elquí|1167",
    "179|44831|2800|2000|2500|2250|Provincia del Elquí|1500",
    "180|44922|2000|2500|3000|2750|Provincia del Elquí|1833",
    "181|44292|3400|2000|2500|2250|Provincia del Elquí|1500",
    "182|44355|3200|1300|1500|1400|Provincia del Elquí|933",
    "183|44987|2200|2000|2500|2250|Provincia del Elquí|1500",
    "184|44572|2900|3000|3500|3250|Provincia del Elquí|2167",
    "185|44245|3200|1500|2000|1750|Provincia del Elquí|1167",
    "186|44518|3320|1300|1500|1400|Provincia del Elquí|933",
    "187|44441|3200|2000|2500|2250|Provincia del Elquí|1500",
    "188|44812|3000|2000|2500|2250|Provincia del Elquí|1500",
    "189|44427|3360|2000|2500|2250|Provincia del Elquí|1500",
    "190|44315|3120|1300|1500|1400|Provincia del Elquí|933",
    "191|45097|2000|1900|2000|1950|Provincia del Elquí|1300",
    "192|44560|3400|2500|3000|2750|Provincia del Elquí|1833",
    "193|44175|3000|1300|1500|1400|Provincia del Elquí|933",
    "194|44952|2000|3000|3500|3250|Provincia del Elquí|2167",
    "195|45090|2000|1800|2000|1900|Provincia del Elquí|1267",
    "196|44413|3360|2000|2500|2250|Provincia del Elquí|1500",
    "197|44784|2600|2000|2500|2250|Provincia del Elquí|1500",
    "198|44224|2800|1300|1500|1400|Provincia del Elquí|933",
    "199|44280|3000|2000|2500|2250|Provincia del Elquí|1500",
    "200|45085|2400|1800|2000|1900|Provincia del Elquí|1267",
    "201|44609|2600|2300|2500|2400|Provincia del Elquí|1600",
    "202|44320|3400|1300|1500|1400|Provincia del Elquí|933",
    "203|44670|3000|2000|2500|2250|Provincia del Elquí|1500",
    "204|45055|2000|2000|2500|2250|Provincia del Elquí|1500",
    "205|45062|2200|2000|2500|2250|Provincia del Elquí|1500",
    "206|45008|2300|1800|2000|1900|Provincia del Elquí|1267",
    "207|44504|3200|1300|1500|1400|Provincia del Elquí|933",
    "208|44838|2800|2000|2500|2250|Provincia del Elquí|1500",
    "209|44719|3200|1500|2000|1750|Provincia del Elquí|1167",
    "210|44677|2400|2500|3000|2750|Provincia del Elquí|1833",
    "211|44742|3000|1300|1500|1400|Provincia del Elquí|933",
    "212|44490|3200|1300|1500|1400|Provincia del Elquí|933",
    "213|45071|2400|1800|2000|1900|Provincia del Elquí|1267",
    "214|45127|2400|2500|3000|2750|Provincia del Elquí|1833",
    "215|44376|3200|1500|2000|1750|Provincia del Elquí|1167",
    "216|44273|3000|2000|2500|2250|Provincia del Elquí|1500",
    "217|44215|2800|1300|1500|1400|Provincia del Elquí|933",
    "218|44945|2400|3000|3500|3250|Provincia del Elquí|2167",
    "219|44714|3200|1500|2000|1750|Provincia del Elquí|1167",
    "220|44826|3000|2000|2500|2250|Provincia del Elquí|1500",
    "221|44565|3000|3000|3500|3250|Provincia del Elquí|2167",
    "222|44901|1600|3000|3500|3250|Provincia del Elquí|2167",
    "223|45132|2000|2500|3000|2750|Provincia del Elquí|1833",
    "224|44322|3320|1300|1500|1400|Provincia del Elquí|933",
    "225|44495|2860|1300|1500|1400|Provincia del Elquí|933",
    "226|45111|2400|2500|3000|2750|Provincia del Elquí|1833",
    "227|44511|3360|1300|1500|1400|Provincia del Elquí|933",
    "228|44420|3400|2000|2500|2250|Provincia del Elquí|1500",
    "229|44924|2000|3000|3500|3250|Provincia del Elquí|2167",
    "230|44971|2000|2500|3000|2750|Provincia del Elquí|1833",
    "231|44364|3200|1500|2000|1750|Provincia del Elquí|1167",
    "232|44985|2500|2000|2500|2250|Provincia del Elquí|1500",
    "233|44644|2400|2300|2500|2400|Provincia del Elquí|1600",
    "234|44630|2000|2500|3000|2750|Provincia del Elquí|1833",
    "235|44859|2400|1500|2000|1750|Provincia de Limarí|1167",
    "236|44637|2460|2500|3000|2750|Provincia del Elquí|1833",
    "237|45112|2400|2500|3000|2750|Provincia del Elquí|1833",
    "238|44348|3360|1300|1500|1400|Provincia del Elquí|933",
    "239|44574|3200|3000|3500|3250|Provincia del Elquí|2167",
    "240|44749|3000|1500|2000|1750|Provincia del Elquí|1167",
    "241|44957|2000|3000|3500|3250|Provincia del Elquí|2167",
    "242|44763|2000|2000|2500|2250|Provincia del Elquí|1500",
    "243|45099|2000|1800|2000|1900|Provincia del Elquí|1267",
    "244|44341|3360|1300|1500|1400|Provincia del Elquí|933",
    "245|44607|2400|2300|2500|2400|Provincia del Elquí|1600"
)

# Date column (D) uses a custom date number format; grab it from an existing
# data row so the brand-new row 245 gets the same cell style/format.
$dateFormat = $ws.Cells.Item(162, 4).NumberFormat

foreach ($line in $rowsData) {
    $parts = $line.Split('|')
    $r = [int]$parts[0]
    $d = [double]$parts[1]
    $j = [double]$parts[2]
    $k = [double]$parts[3]
    $l = [double]$parts[4]
    $m = [double]$parts[5]
    $o = $parts[6]
    $p = [double]$parts[7]

    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
}

# New row 245 needs the remaining (previously constant / carried-over) columns too,
# matching the rest of the "Perejil" rows exactly.
$carryCols = @(1, 2, 3, 5, 6, 7, 8, 9, 14, 17, 18)
foreach ($c in $carryCols) {
    $srcVal = $ws.Cells.Item(244, $c).Value()
    $ws.Cells.Item(245, $c).Value = $srcVal
}

Write-Output "Applied weekly shift to rows 161-245"
